$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("K5").Value = 19.30324074074072
$ws.Range("R5").Value = 1.975772235794973
$ws.Range("S5").Value = 2.177153507468733
$ws.Range("K6").Value = 19.30324074074072
$ws.Range("R6").Value = 1.352319749654237
$ws.Range("S6").Value = 1.423118895050623
$ws.Range("K7").Value = 19.30324074074072
$ws.Range("K11").Value = 12.67039049919483
$ws.Range("R11").Value = 1.847705673092716
$ws.Range("S11").Value = 2.020749367497032
$ws.Range("K12").Value = 12.67039049919483
$ws.Range("R12").Value = 1.300269876134972
$ws.Range("S12").Value = 1.364969709970879
$ws.Range("K13").Value = 12.67039049919483
$ws.Range("K14").Value = 13.17361111111111
$ws.Range("K15").Value = 13.17361111111111
$ws.Range("R15").Value = 1.304077921028169
$ws.Range("S15").Value = 1.369214264257821
$ws.Range("K16").Value = 13.17361111111111
$ws.Range("R16").Value = 1.856836936506854
$ws.Range("S16").Value = 2.031823338122968
$ws.Range("K17").Value = 12.92654320987656
$ws.Range("R17").Value = 1.302205489329493
$ws.Range("S17").Value = 1.367127007643996
$ws.Range("K18").Value = 12.92654320987656
$ws.Range("K19").Value = 12.92654320987656
$ws.Range("R19").Value = 1.852342479090948
$ws.Range("S19").Value = 2.026371187792892
$ws.Range("K20").Value = -1.226851851851833
$ws.Range("R20").Value = 1.626775542720574
$ws.Range("S20").Value = 1.756382654173023
$ws.Range("K21").Value = -1.226851851851833
$ws.Range("K22").Value = -1.226851851851833
$ws.Range("R22").Value = 1.203236793039155
$ws.Range("S22").Value = 1.257328254301852
$ws.Range("K23").Value = 13.0158303464755
$ws.Range("R23").Value = 1.302881541082627
$ws.Range("S23").Value = 1.367880580392128
$ws.Range("K24").Value = 13.0158303464755
$ws.Range("K25").Value = 13.0158303464755
$ws.Range("R25").Value = 1.853964204859962
$ws.Range("S25").Value = 2.02833814451736
$ws.Range("K29").Value = 1.925925925925943
$ws.Range("K30").Value = 1.925925925925943
$ws.Range("R30").Value = 1.672133966590239
$ws.Range("S30").Value = 1.810106000718649
$ws.Range("K31").Value = 1.925925925925943
$ws.Range("R31").Value = 1.223958122597613
$ws.Range("S31").Value = 1.280232184891932
$ws.Range("K32").Value = -1.819444444444444
$ws.Range("K33").Value = -1.819444444444444
$ws.Range("R33").Value = 1.199420117463385
$ws.Range("S33").Value = 1.253114445055376
$ws.Range("K34").Value = -1.819444444444444
$ws.Range("R34").Value = 1.618523362263702
$ws.Range("S34").Value = 1.746638928617865
$ws.Range("K35").Value = 13.0158303464755
$ws.Range("R35").Value = 1.302881541082627
$ws.Range("S35").Value = 1.367880580392128
$ws.Range("K36").Value = 13.0158303464755
$ws.Range("R36").Value = 1.853964204859962
$ws.Range("S36").Value = 2.02833814451736
$ws.Range("K37").Value = 13.0158303464755
$ws.Range("K38").Value = 5.486111111111112
$ws.Range("K39").Value = 5.486111111111112
$ws.Range("R39").Value = 1.248232108317215
$ws.Range("S39").Value = 1.307120032773454
$ws.Range("K40").Value = 5.486111111111112
$ws.Range("R40").Value = 1.726493341788205
$ws.Range("S40").Value = 1.874863921842289
$ws.Range("K44").Value = -1.226851851851833
$ws.Range("K45").Value = -1.226851851851833
$ws.Range("R45").Value = 1.626775542720574
$ws.Range("S45").Value = 1.756382654173023
$ws.Range("K46").Value = -1.226851851851833
$ws.Range("R46").Value = 1.203236793039155
$ws.Range("S46").Value = 1.257328254301852
$ws.Range("K47").Value = 3.38888888888889
$ws.Range("R47").Value = 1.233817681248088
$ws.Range("S47").Value = 1.291146001942376
$ws.Range("K48").Value = 3.38888888888889
$ws.Range("R48").Value = 1.694051767048283
$ws.Range("S48").Value = 1.836167304537999
$ws.Range("K49").Value = 3.38888888888889
$ws.Range("K56").Value = 16.86342592592595
$ws.Range("R56").Value = 1.92665172779809
$ws.Range("S56").Value = 2.116885095206829
$ws.Range("K57").Value = 16.86342592592595
$ws.Range("R57").Value = 1.332696358504853
$ws.Range("S57").Value = 1.401162263046183
$ws.Range("K58").Value = 16.86342592592595
